# Chap04: Correction Lucien finished
#
# 1) The "Fixed date" footer placeholder (datetimeFigureOut field) on the
#    slide master and on every slide layout is corrected from 31/10/2017
#    to 09/11/2017.
# 2) A handful of pictures/shapes on slide 1 are nudged vertically (and one
#    is also resized) to re-align the figure after the date correction.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fix the cached date text wherever it appears (master + every layout).
# ---------------------------------------------------------------------------
function Update-DateFooter($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.HasTextFrame) {
      $tf = $sh.TextFrame
      if ($tf.HasText) {
        if ($tf.TextRange.Text -eq "31/10/2017") {
          $tf.TextRange.Text = "09/11/2017"
        }
      }
    }
  }
}

$master = $p.SlideMaster
Update-DateFooter $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
  $layout = $master.CustomLayouts.Item($L)
  Update-DateFooter $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Re-position/resize the pictures and call-out shapes on slide 1.
#    Top/Height are expressed in points (1 pt = 12700 EMU); the literals
#    below are chosen so that the underlying EMU values match exactly.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# Picture 2 -> y: -15169 EMU -> -53751 EMU
$slide.Shapes.Item(1).Top = -4.232401371002197

# Picture 3 -> y: 5169402 EMU -> 5130820 EMU
$slide.Shapes.Item(2).Top = 404.0016174316406

# Groupe 109 -> y: 200854 EMU -> 162272 EMU
$slide.Shapes.Item(3).Top = 12.777362823486328

# ZoneTexte 156 "(a)" -> y: 200855 EMU -> 162273 EMU
$slide.Shapes.Item(4).Top = 12.777441024780273

# ZoneTexte 157 "(c)" -> y: 129095 EMU -> 166618 EMU
$slide.Shapes.Item(5).Top = 13.119566917419434

# ZoneTexte 158 "(b)" -> y: 5457439 EMU -> 5418857 EMU
$slide.Shapes.Item(6).Top = 426.68170166015625

# Rectangle 159 -> y: 56839 EMU -> 195699 EMU ; cy: 9073008 EMU -> 8895566 EMU
$slide.Shapes.Item(7).Top = 15.409409523010254
$slide.Shapes.Item(7).Height = 700.4382934570312
